$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.544.91"
$ws.Range("E2").Value = "  -0.88%  "

$ws.Range("D3").Value = "2.222.17"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.34"
$ws.Range("E5").Value = "  +7.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.627"
$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.69"
$ws.Range("E7").Value = "  +2.55%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.580"
$ws.Range("E9").Value = "  +5.17%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.26"
$ws.Range("E10").Value = "  +17.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0962"
$ws.Range("E11").Value = "  -1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.41"
$ws.Range("E12").Value = "  +0.83%  "

$ws.Range("E13").Value = "  +1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +3.77%  "

$ws.Range("D15").Value = "2.552.87"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.84"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.854"
$ws.Range("E17").Value = "  -0.03%  "

$ws.Range("D18").Value = "2.218.75"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").Value = "41.511.32"
$ws.Range("E19").Value = "  -0.79%  "

$ws.Range("D20").Value = "0.0₃0966"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.19"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.48"
$ws.Range("E22").Value = "  -0.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.85"
$ws.Range("E23").Value = "  -0.79%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.23"
$ws.Range("E24").Value = "  +10.42%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.89"
$ws.Range("E25").Value = "  +6.94%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.49"
$ws.Range("E27").Value = "  +6.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.44"
$ws.Range("E28").Value = "  +4.90%  "

$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.24"
$ws.Range("E30").Value = "  +1.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.58"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("E32").Value = "  +1.87%  "

$ws.Range("E33").Value = "  -1.94%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.53"
$ws.Range("E34").Value = "  +3.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.67"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.88"
$ws.Range("E37").Value = "  +18.32%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.91"
$ws.Range("E38").Value = "  +9.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0292"
$ws.Range("E39").Value = "  +10.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  +1.68%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.60"
$ws.Range("E41").Value = "  +3.91%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.92"
$ws.Range("E42").Value = "  -1.06%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.208"
$ws.Range("E43").Value = "  +9.89%  "

$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.75"
$ws.Range("E44").Value = "  +18.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.75"
$ws.Range("E46").Value = "  -2.02%  "

$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.75"
$ws.Range("E47").Value = "  +11.34%  "

$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("E50").Value = "  +6.95%  "

$ws.Range("E51").Value = "  +1.68%  "
